$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: AK1 = "BOUNDARY" header (merged AK1:AR1), AL1:AR1 stay blank but styled ---
$ws.Range("F1").Copy($ws.Range("AK1:AR1"))
$ws.Range("AK1").Value = "BOUNDARY"

# --- Row 2: epsilon labels as text, styled like existing header row (bold/border/center) ---
$scratch = $ws.Range("ZZ1")
$ws.Range("AJ2").Copy($scratch)
$scratch.NumberFormat = "@"
$scratch.Copy($ws.Range("AK2"))
$scratch.Copy($ws.Range("AL2"))
$scratch.Copy($ws.Range("AM2"))
$scratch.Copy($ws.Range("AN2"))
$scratch.Copy($ws.Range("AO2"))
$scratch.Copy($ws.Range("AP2"))
$scratch.Copy($ws.Range("AQ2"))
$scratch.Copy($ws.Range("AR2"))
$scratch.Clear()
$ws.Range("AK2").Value = "0.01"
$ws.Range("AL2").Value = "0.02"
$ws.Range("AM2").Value = "0.03"
$ws.Range("AN2").Value = "0.04"
$ws.Range("AO2").Value = "0.05"
$ws.Range("AP2").Value = "0.07"
$ws.Range("AQ2").Value = "0.10"
$ws.Range("AR2").Value = "0.20"

# --- Row 4: BOUNDARY attack numeric results ---
$ws.Range("AK4").Value = 432.7917705790202
$ws.Range("AL4").Value = 435.5175196075439
$ws.Range("AM4").Value = 439.2190765889486
$ws.Range("AN4").Value = 442.2330904579163
$ws.Range("AO4").Value = 447.3975617790222
$ws.Range("AP4").Value = 462.42976770401
$ws.Range("AQ4").Value = 474.4381432723999
$ws.Range("AR4").Value = 598.2050880622863

# --- Row 5: BOUNDARY attack numeric results ---
$ws.Range("AK5").Value = 540.9919002676994
$ws.Range("AL5").Value = 544.8139485716558
$ws.Range("AM5").Value = 548.5361284987341
$ws.Range("AN5").Value = 549.0153975015104
$ws.Range("AO5").Value = 555.8957974614127
$ws.Range("AP5").Value = 581.4198112572533
$ws.Range("AQ5").Value = 596.6181559478788
$ws.Range("AR5").Value = 783.0152101056404

# --- Row 6: BOUNDARY attack numeric results ---
$ws.Range("AK6").Value = 0.9990363688847373
$ws.Range("AL6").Value = 0.9990245789996384
$ws.Range("AM6").Value = 0.9990095397240493
$ws.Range("AN6").Value = 0.9990069293620361
$ws.Range("AO6").Value = 0.9989863542368618
$ws.Range("AP6").Value = 0.9988941329854391
$ws.Range("AQ6").Value = 0.9988331431111914
$ws.Range("AR6").Value = 0.9980908408034129

# --- Row 7: BOUNDARY attack numeric results ---
$ws.Range("AK7").Value = 399.0320683479309
$ws.Range("AL7").Value = 415.4160107549031
$ws.Range("AM7").Value = 405.1621262041728
$ws.Range("AN7").Value = 417.024144077301
$ws.Range("AO7").Value = 432.9904280662537
$ws.Range("AP7").Value = 444.4669622866313
$ws.Range("AQ7").Value = 544.9597953160604
$ws.Range("AR7").Value = 700.3300581804912

# --- Row 8: BOUNDARY attack numeric results ---
$ws.Range("AK8").Value = 520.4198570710514
$ws.Range("AL8").Value = 540.0463791846861
$ws.Range("AM8").Value = 525.0756877510274
$ws.Range("AN8").Value = 538.6864789263687
$ws.Range("AO8").Value = 549.551265392308
$ws.Range("AP8").Value = 568.514633579456
$ws.Range("AQ8").Value = 678.2156310702549
$ws.Range("AR8").Value = 863.5437383340842

# --- Row 9: BOUNDARY attack numeric results ---
$ws.Range("AK9").Value = 0.9992002936530356
$ws.Range("AL9").Value = 0.9991314404105985
$ws.Range("AM9").Value = 0.9991916481528849
$ws.Range("AN9").Value = 0.9991407673813675
$ws.Range("AO9").Value = 0.9991265007440792
$ws.Range("AP9").Value = 0.9990315614699359
$ws.Range("AQ9").Value = 0.9984186432039401
$ws.Range("AR9").Value = 0.9973475548268412

# --- Row 10: BOUNDARY attack numeric results ---
$ws.Range("AK10").Value = 283.2202223014831
$ws.Range("AL10").Value = 296.9328125190735
$ws.Range("AM10").Value = 305.4126160240173
$ws.Range("AN10").Value = 306.4555748875936
$ws.Range("AO10").Value = 314.624032535553
$ws.Range("AP10").Value = 341.1765812110901
$ws.Range("AQ10").Value = 410.6481328964234
$ws.Range("AR10").Value = 666.5929933834076

# --- Row 11: BOUNDARY attack numeric results ---
$ws.Range("AK11").Value = 402.7637699362363
$ws.Range("AL11").Value = 418.2314816654203
$ws.Range("AM11").Value = 416.5085486369327
$ws.Range("AN11").Value = 421.971121361415
$ws.Range("AO11").Value = 429.1832709963401
$ws.Range("AP11").Value = 459.8072883257767
$ws.Range("AQ11").Value = 517.8801462060561
$ws.Range("AR11").Value = 840.7632142560259

# --- Row 12: BOUNDARY attack numeric results ---
$ws.Range("AK12").Value = 0.9993913910429063
$ws.Range("AL12").Value = 0.9993391793864905
$ws.Range("AM12").Value = 0.999351846827865
$ws.Range("AN12").Value = 0.99932198852751
$ws.Range("AO12").Value = 0.999297560386176
$ws.Range("AP12").Value = 0.9991944463183169
$ws.Range("AQ12").Value = 0.998978395251458
$ws.Range("AR12").Value = 0.9973350900343017

# --- Merge the new BOUNDARY header range ---
$ws.Range("AK1:AR1").Merge()

Write-Output "done"